$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose highlighted (red-fill) style is removed along with a new computed value
$ws.Range("C2").Value = 17.40819215506482
$ws.Range("C2").Style = "Normal"
$ws.Range("I2").Value = 0.02057613168724279
$ws.Range("I2").Style = "Normal"
$ws.Range("C3").Value = 19.91810114448176
$ws.Range("C3").Style = "Normal"
$ws.Range("I3").Value = 0.03645833333333334
$ws.Range("I3").Style = "Normal"
$ws.Range("C4").Value = 19.12800902518529
$ws.Range("C4").Style = "Normal"
$ws.Range("I4").Value = 0.03645833333333334
$ws.Range("I4").Style = "Normal"
$ws.Range("C5").Value = 15.74855722832491
$ws.Range("C5").Style = "Normal"
$ws.Range("I5").Value = 0.05333333333333334
$ws.Range("I5").Style = "Normal"
$ws.Range("K5").Value = 3.485826843761838
$ws.Range("K5").Style = "Normal"
$ws.Range("C6").Value = 19.60514228209162
$ws.Range("C6").Style = "Normal"
$ws.Range("I6").Value = 0.03645833333333334
$ws.Range("I6").Style = "Normal"
$ws.Range("C7").Value = 18.25902148010383
$ws.Range("C7").Style = "Normal"
$ws.Range("I7").Value = 0.03645833333333334
$ws.Range("I7").Style = "Normal"
$ws.Range("C8").Value = 15.07024562620989
$ws.Range("C8").Style = "Normal"
$ws.Range("I8").Value = 0.05333333333333334
$ws.Range("I8").Style = "Normal"
$ws.Range("K8").Value = 3.205331201526457
$ws.Range("K8").Style = "Normal"
$ws.Range("C9").Value = 19.87341873939329
$ws.Range("C9").Style = "Normal"
$ws.Range("I9").Value = 0.03645833333333334
$ws.Range("I9").Style = "Normal"
$ws.Range("C10").Value = 19.8912023015621
$ws.Range("C10").Style = "Normal"
$ws.Range("I10").Value = 0.03645833333333334
$ws.Range("I10").Style = "Normal"
$ws.Range("C11").Value = 13.95569405639475
$ws.Range("C11").Style = "Normal"
$ws.Range("I11").Value = 0.03856749311294766
$ws.Range("I11").Style = "Normal"
$ws.Range("K11").Value = 3.447937982072726
$ws.Range("K11").Style = "Normal"
$ws.Range("C12").Value = 19.31401097013146
$ws.Range("C12").Style = "Normal"
$ws.Range("I12").Value = 0.03645833333333334
$ws.Range("I12").Style = "Normal"
$ws.Range("C13").Value = 18.77543457306755
$ws.Range("C13").Style = "Normal"
$ws.Range("I13").Value = 0.005208333333333334
$ws.Range("I13").Style = "Normal"
$ws.Range("C14").Value = 14.71971926437311
$ws.Range("C14").Style = "Normal"
$ws.Range("I14").Value = 0.05333333333333334
$ws.Range("I14").Style = "Normal"
$ws.Range("K14").Value = 2.763424943782517
$ws.Range("K14").Style = "Normal"
$ws.Range("C15").Value = 19.47582887855932
$ws.Range("C15").Style = "Normal"
$ws.Range("I15").Value = 0.03645833333333334
$ws.Range("I15").Style = "Normal"
$ws.Range("C16").Value = 15.05564418777329
$ws.Range("C16").Style = "Normal"
$ws.Range("I16").Value = 0.02057613168724279
$ws.Range("I16").Style = "Normal"
$ws.Range("C17").Value = 13.60061437498968
$ws.Range("C17").Style = "Normal"
$ws.Range("I17").Value = 0.04683195592286502
$ws.Range("I17").Style = "Normal"
$ws.Range("K17").Value = 2.269837011114737
$ws.Range("K17").Style = "Normal"
$ws.Range("M17").Value = 1
$ws.Range("M17").Style = "Normal"
$ws.Range("C18").Value = 18.13592238629358
$ws.Range("C18").Style = "Normal"
$ws.Range("I18").Value = 0.02083333333333333
$ws.Range("I18").Style = "Normal"
$ws.Range("C19").Value = 20.3637157720036
$ws.Range("C19").Style = "Normal"
$ws.Range("I19").Value = 0.03645833333333334
$ws.Range("I19").Style = "Normal"
$ws.Range("C20").Value = 15.62873760559645
$ws.Range("C20").Style = "Normal"
$ws.Range("I20").Value = 0.05333333333333334
$ws.Range("I20").Style = "Normal"
$ws.Range("K20").Value = 2.983229197536327
$ws.Range("K20").Style = "Normal"
$ws.Range("C21").Value = 19.53596071976449
$ws.Range("C21").Style = "Normal"
$ws.Range("I21").Value = 0.03645833333333334
$ws.Range("I21").Style = "Normal"
$ws.Range("C22").Value = 17.21209395200373
$ws.Range("C22").Style = "Normal"
$ws.Range("I22").Value = 0.02057613168724279
$ws.Range("I22").Style = "Normal"
$ws.Range("C23").Value = 16.82955425677702
$ws.Range("C23").Style = "Normal"
$ws.Range("I23").Value = 0.02083333333333333
$ws.Range("I23").Style = "Normal"
$ws.Range("C24").Value = 19.11706532842858
$ws.Range("C24").Style = "Normal"
$ws.Range("I24").Value = 0.02083333333333333
$ws.Range("I24").Style = "Normal"
$ws.Range("C25").Value = 12.81366767882287
$ws.Range("C25").Style = "Normal"
$ws.Range("I25").Value = 0.01333333333333333
$ws.Range("I25").Style = "Normal"
$ws.Range("K25").Value = 2.368434805890583
$ws.Range("K25").Style = "Normal"
$ws.Range("M25").Value = 0.4583333333333333
$ws.Range("M25").Style = "Normal"
$ws.Range("C26").Value = 14.78490457790944
$ws.Range("C26").Style = "Normal"
$ws.Range("I26").Value = 0.02333333333333333
$ws.Range("I26").Style = "Normal"
$ws.Range("K26").Value = 3.519394837090817
$ws.Range("K26").Style = "Normal"
$ws.Range("C27").Value = 20.0659061787455
$ws.Range("C27").Style = "Normal"
$ws.Range("I27").Value = 0.03645833333333334
$ws.Range("I27").Style = "Normal"
$ws.Range("C28").Value = 13.65491720725454
$ws.Range("C28").Style = "Normal"
$ws.Range("I28").Value = 0.04333333333333335
$ws.Range("I28").Style = "Normal"
$ws.Range("K28").Value = 3.707835569091639
$ws.Range("K28").Style = "Normal"
$ws.Range("C29").Value = 18.84136784189097
$ws.Range("C29").Style = "Normal"
$ws.Range("I29").Value = 0.03645833333333334
$ws.Range("I29").Style = "Normal"
$ws.Range("C30").Value = 18.34457096166117
$ws.Range("C30").Style = "Normal"
$ws.Range("I30").Value = 0.03645833333333334
$ws.Range("I30").Style = "Normal"
$ws.Range("C31").Value = 18.92704318051282
$ws.Range("C31").Style = "Normal"
$ws.Range("I31").Value = 0.02083333333333333
$ws.Range("I31").Style = "Normal"
$ws.Range("C32").Value = 16.50390303992081
$ws.Range("C32").Style = "Normal"
$ws.Range("I32").Value = 0.02057613168724279
$ws.Range("I32").Style = "Normal"
$ws.Range("C33").Value = 19.37507748769768
$ws.Range("C33").Style = "Normal"
$ws.Range("I33").Value = 0.03645833333333334
$ws.Range("I33").Style = "Normal"
$ws.Range("C34").Value = 19.82103789235044
$ws.Range("C34").Style = "Normal"
$ws.Range("I34").Value = 0.03645833333333334
$ws.Range("I34").Style = "Normal"
$ws.Range("C35").Value = 15.01679477297077
$ws.Range("C35").Style = "Normal"
$ws.Range("I35").Value = 0.02333333333333334
$ws.Range("I35").Style = "Normal"
$ws.Range("K35").Value = 3.096818523806561
$ws.Range("K35").Style = "Normal"
$ws.Range("C36").Value = 18.95956737585637
$ws.Range("C36").Style = "Normal"
$ws.Range("I36").Value = 0.005208333333333334
$ws.Range("I36").Style = "Normal"
$ws.Range("C37").Value = 19.6978788678873
$ws.Range("C37").Style = "Normal"
$ws.Range("I37").Value = 0.03645833333333334
$ws.Range("I37").Style = "Normal"

# Cells that keep their existing style (s=1 red fill) but get an updated value
$ws.Range("M2").Value = 1.076923076923077
$ws.Range("M3").Value = 1.066666666666667
$ws.Range("M4").Value = 1
$ws.Range("M5").Value = 1.041666666666667
$ws.Range("M6").Value = 1.137931034482759
$ws.Range("M7").Value = 1.137931034482759
$ws.Range("M8").Value = 1.166666666666667
$ws.Range("M9").Value = 1.206896551724138
$ws.Range("M10").Value = 1.033333333333333
$ws.Range("M11").Value = 1.045454545454545
$ws.Range("M12").Value = 1.033333333333333
$ws.Range("M13").Value = 1
$ws.Range("M14").Value = 0.9583333333333334
$ws.Range("M15").Value = 1.066666666666667
$ws.Range("M16").Value = 0.6153846153846154
$ws.Range("M18").Value = 0.896551724137931
$ws.Range("M19").Value = 1.033333333333333
$ws.Range("M20").Value = 0.9166666666666666
$ws.Range("M21").Value = 1.033333333333333
$ws.Range("M22").Value = 1.038461538461539
$ws.Range("M23").Value = 0.5357142857142857
$ws.Range("M24").Value = 1.142857142857143
$ws.Range("M26").Value = 1.08695652173913
$ws.Range("M27").Value = 1.137931034482759
$ws.Range("M28").Value = 1.260869565217391
$ws.Range("M29").Value = 1.172413793103448
$ws.Range("M30").Value = 1.068965517241379
$ws.Range("M31").Value = 1.03448275862069
$ws.Range("M32").Value = 1.076923076923077
$ws.Range("M33").Value = 1.066666666666667
$ws.Range("M34").Value = 1.03448275862069
$ws.Range("M35").Value = 1
$ws.Range("M36").Value = 1.137931034482759
$ws.Range("M37").Value = 1
